# "Fixed intents and schema"
# Adds three new smartphoneName rows ("Redmi Note 4", "P20", "P20 lite") to the
# entities table, right before the existing "smartphoneRange" block (old row 21),
# and folds them into the existing "smartphoneName" merged group (A12:A20 -> A12:A23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the old row 21 ("smartphoneRange" / "Gama baja" ...),
# pushing everything below (old rows 21-34) down to rows 24-37.
$ws.Rows("21:23").Insert()

# New synonym entries for the "smartphoneName" entity.
$ws.Range("B21").Value = "Redmi Note 4"
$ws.Range("C21").Value = "Note 4" + [char]10 + "redmi 4" + [char]10 + "redmi note 4"

$ws.Range("B22").Value = "P20"
$ws.Range("C22").Value = "p 20"

$ws.Range("B23").Value = "P20 lite"
$ws.Range("C23").Value = "p 20 lite"

# Match the rest of the table's formatting: column A/B centered, no wrap;
# column C centered with wrapped multi-line synonyms.
$ws.Range("A21:B23").HorizontalAlignment = -4108
$ws.Range("A21:B23").VerticalAlignment = -4108
$ws.Range("A21:B23").WrapText = $false

$ws.Range("C21:C23").HorizontalAlignment = -4108
$ws.Range("C21:C23").VerticalAlignment = -4108
$ws.Range("C21:C23").WrapText = $true

# Row 21 now holds 3 lines of synonyms, so it needs extra height.
$ws.Rows(21).RowHeight = 45

# Extend the "smartphoneName" merged header cell to span the new rows too.
$ws.Range("A12:A23").Merge() | Out-Null

# Leave the selection where the edit was made.
$ws.Range("B21").Select() | Out-Null
